$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.576.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.604.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.22%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.988"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.517"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.984"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.60%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.03"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.254"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0603"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0869"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.842.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.63%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.615.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.554"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.623.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "226.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.90%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0721"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.989"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.61%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.110"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.989"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.64%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0476"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.370.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.78%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.985"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0167"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.83%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.548"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.842"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.61%  "

$ws.Range("E41").Value = "  -1.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.987"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.34%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.27%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.21%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.97%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.764.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.96%  "

$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.100"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.25%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0500"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.23%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.47%  "
